$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 203, shifting existing rows 203:239 down to 204:240
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new data record
$ws.Cells.Item(203, 1).Value = 3
$ws.Cells.Item(203, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(203, 3).Value = "Coquimbo"
$ws.Cells.Item(203, 4).Value = 44522
$ws.Cells.Item(203, 5).Value = 5
$ws.Cells.Item(203, 6).Value = 100114013
$ws.Cells.Item(203, 7).Value = "Zanahoria"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 310
$ws.Cells.Item(203, 11).Value = 6500
$ws.Cells.Item(203, 12).Value = 7000
$ws.Cells.Item(203, 13).Value = 6742
$ws.Cells.Item(203, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(203, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(203, 16).Value = 337
$ws.Cells.Item(203, 17).Value = 20
$ws.Cells.Item(203, 18).Value = "Hortaliza"
